# Results_per_Country 2050_BG.xlsx — corrected-code results update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 (Hydrogen / Non-metallic minerals): corrected value removed -> blank cell
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = $ws.Range("C3").Style

# C4 (Methanol / Chemicals): value corrected to 0
$ws.Range("C4").Value = 0

# C5 (Ammonia / Chemicals): value corrected
$ws.Range("C5").Value = 1219.956671505592

# Row 7: "Other" is renamed to "Biogas" and its value corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 2367.708965642697

# New row 8 ("Other") appears below the renamed Biogas row, with its own value
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"

$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = $ws.Range("B7").Style

$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = $ws.Range("C7").Style

$ws.Range("D8").Value = 1785.796439581564
